$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the front (A and B) so the existing
# Method / Average Significant Percentage data shifts right to C/D.
$ws.Range("A:B").Insert()

# New header cells (Network, Alpha) need the same bold/bordered/centered
# style already used by the other header cells - copy it over.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$ws.Range("A1").Value = "Network"
$ws.Range("B1").Value = "Alpha"

# Fill the Network / Alpha columns for each data row.
$ws.Range("A2:A5").Value = "HumanNet"
$ws.Range("B2:B5").Value = 0.1

# Update the Average Significant Percentage values (now column D) to
# the new alpha=0.1 results.
$ws.Range("D2").Value = 1.732068642831585
$ws.Range("D3").Value = 0.02270663033605813
$ws.Range("D4").Value = 0.04668620813034981
$ws.Range("D5").Value = 0
